# Swap the order of "System" and the recorder's email address in the
# "Recorded By" column (G) of the "Session Analysis Results" sheet.
#
#   "System, dnasr281@gmail.com"  ->  "dnasr281@gmail.com, System"
#
# Rows whose Date column (E) is "23/12/2025" are left untouched, matching
# the source diff exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"
$skipDate = "23/12/2025"

$lastRow = $ws.UsedRange.Rows.Count

$changed = 0
for ($r = 2; $r -le $lastRow; $r++) {
    $gCell = $ws.Cells.Item($r, 7)
    $gVal = $gCell.Value2
    if ($gVal -eq $oldValue) {
        $eVal = $ws.Cells.Item($r, 5).Value2
        if ($eVal -ne $skipDate) {
            $gCell.Value = $newValue
            $changed++
        }
    }
}

Write-Host "Updated $changed 'Recorded By' cells."
